# Applies the "Add files via upload" commit:
#  - swap the B10/C10 values
#  - fill in the previously-empty rows 11-19 with the Ano/M1/Precos data
#  - add a brand-new row 20 (with a D20 cell using the default style)
#  - give C20 a new "#,##0.0" centered number format (numFmt 166 / new cellXf)
#  - nudge column A to a custom width and move the active selection to C20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: values were swapped (B10 now holds 15, C10 now holds 16.7) ---
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 16.7

# --- Rows 11-19: previously blank, now populated ---
$yearData = @(
  @(1946, 10.6, 16.5),
  @(1947, -0.9, 21.9),
  @(1948, 6.7, 3.4),
  @(1949, 17.1, 4.3),
  @(1950, 27.5, 9.4),
  @(1951, 18.2, 12.1),
  @(1952, 16.8, 17.3),
  @(1953, 19.6, 14.3),
  @(1954, 22.4, 22.6)
)

$r = 11
foreach ($row in $yearData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Row 20: new row, including D20 which keeps the sheet's default style ---
$ws.Range("A20").Value = 1955
$ws.Range("B20").Value = 19.4
$ws.Range("C20").Value = 23
$ws.Range("D20").Value = 2

# C20 gets its own number format (adds numFmt 166 "#,##0.0" + a new centered cellXf)
$ws.Range("C20").NumberFormat = "#,##0.0"
$ws.Range("C20").HorizontalAlignment = -4108

# --- Column A: mark the width as explicit/custom (was using the sheet default) ---
$ws.Columns.Item(1).ColumnWidth = 10.6

# --- Move the selection/active cell to C20 ---
$ws.Range("C20").Select()
